$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLogin")

# --- Insert a new test-case row at row 5 -----------------------------------
# (pushes the existing rows 5-10 down to 6-11; Excel copies the row's
# formatting down automatically, so the new row keeps the same look as the
# rest of the table)
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "standard_user@"
$ws.Range("B5").Value = "secret_sauce"
$ws.Range("C5").Value = "Epic sadface: Username and password do not match any user in this service"

# --- Append two more test-case rows at the bottom --------------------------
# Inserting (rather than just typing past the end) keeps the same formatting
# as the row above.
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "standard_user"
$ws.Range("B12").Value = "secret_sauce#"
$ws.Range("C12").Value = "Epic sadface: Username and password do not match any user in this service"

$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "locked_out_user"
$ws.Range("B13").Value = "secret_sauce"
$ws.Range("C13").Value = "Epic sadface: Sorry, this user has been locked out."
